# Auto-generated Excel COM-interop script
# Applies market-price / profit-value updates scraped by the scheduled runner
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 292.93103
$ws.Range("I33").Value = 194.82608
$ws.Range("K33").Value = 194.82608
$ws.Range("M33").Value = 34.17392000000001
$ws.Range("H43").Value = 719.6923
$ws.Range("I43").Value = 283.8
$ws.Range("J43").Value = 992.125
$ws.Range("K43").Value = 283.8
$ws.Range("L43").Value = 992.125
$ws.Range("M43").Value = -214.8
$ws.Range("N43").Value = -1130.125
$ws.Range("H98").Value = 1826
$ws.Range("I98").Value = 1968.1666
$ws.Range("J98").Value = 1399.5
$ws.Range("K98").Value = 1968.1666
$ws.Range("L98").Value = 1399.5
$ws.Range("M98").Value = -470.1666
$ws.Range("N98").Value = -4395.5
$ws.Range("H112").Value = 1470.4286
$ws.Range("I112").Value = 799.6667
$ws.Range("J112").Value = 1582.2222
$ws.Range("K112").Value = 2399.0001
$ws.Range("L112").Value = 4746.6666
$ws.Range("M112").Value = -1291.0001
$ws.Range("N112").Value = -6962.6666
$ws.Range("H122").Value = 1826
$ws.Range("I122").Value = 1968.1666
$ws.Range("J122").Value = 1399.5
$ws.Range("K122").Value = 5904.4998
$ws.Range("L122").Value = 4198.5
$ws.Range("M122").Value = -3454.4998
$ws.Range("N122").Value = -9098.5
$ws.Range("H138").Value = 4962.26
$ws.Range("I138").Value = 1273.2222
$ws.Range("J138").Value = 6954.34
$ws.Range("K138").Value = 3819.6666
$ws.Range("L138").Value = 20863.02
$ws.Range("M138").Value = 1320.3334
$ws.Range("N138").Value = -31143.02

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20338.396
$ws.Range("I32").Value = 22927.674
$ws.Range("J32").Value = 11275.929
$ws.Range("K32").Value = 22927.674
$ws.Range("L32").Value = 11275.929
$ws.Range("M32").Value = -22640.674
$ws.Range("N32").Value = -11849.929
$ws.Range("H133").Value = 32761
$ws.Range("J133").Value = 32761
$ws.Range("L133").Value = 32761
$ws.Range("N133").Value = -37821

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2387.5757
$ws.Range("I134").Value = 2292.1538
$ws.Range("J134").Value = 2742
$ws.Range("K134").Value = 6876.4614
$ws.Range("L134").Value = 8226
$ws.Range("M134").Value = -4341.4614
$ws.Range("N134").Value = -13296

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 28000
$ws.Range("J59").Value = 28000
$ws.Range("L59").Value = 28000
$ws.Range("N59").Value = -30290

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 49.533333
$ws.Range("J2").Value = 54.076923
$ws.Range("L2").Value = 324.461538
$ws.Range("N2").Value = -550.461538
$ws.Range("H4").Value = 200.5
$ws.Range("I4").Value = 200.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 601.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -489.5
$ws.Range("N4").ClearContents()
$ws.Range("H17").Value = 1722.2222
$ws.Range("J17").Value = 1687.5
$ws.Range("L17").Value = 5062.5
$ws.Range("N17").Value = -5400.5
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 900
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 2700
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -2473
$ws.Range("N20").Value = -9454
$ws.Range("H34").Value = 2391.05
$ws.Range("I34").Value = 211
$ws.Range("J34").Value = 3117.7334
$ws.Range("K34").Value = 633
$ws.Range("L34").Value = 9353.200199999999
$ws.Range("M34").Value = -549
$ws.Range("N34").Value = -9521.200199999999
$ws.Range("H36").Value = 516.3333
$ws.Range("I36").Value = 699.5
$ws.Range("K36").Value = 2098.5
$ws.Range("M36").Value = -1929.5
$ws.Range("H46").Value = 2970.5278
$ws.Range("J46").Value = 2970.5278
$ws.Range("L46").Value = 8911.5834
$ws.Range("N46").Value = -9093.5834
$ws.Range("H58").Value = 3131.818
$ws.Range("J58").Value = 3131.818
$ws.Range("L58").Value = 9395.454000000002
$ws.Range("N58").Value = -9651.454000000002

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 250005020
$ws.Range("I11").Value = 500000030
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 500000030
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -499999891
$ws.Range("N11").Value = -10278
$ws.Range("H70").Value = 5680.174
$ws.Range("I70").Value = 5164.9
$ws.Range("J70").Value = 6076.5386
$ws.Range("K70").Value = 5164.9
$ws.Range("L70").Value = 6076.5386
$ws.Range("M70").Value = -4894.9
$ws.Range("N70").Value = -6616.5386
$ws.Range("H73").Value = 5680.174
$ws.Range("I73").Value = 5164.9
$ws.Range("J73").Value = 6076.5386
$ws.Range("K73").Value = 5164.9
$ws.Range("L73").Value = 6076.5386
$ws.Range("M73").Value = -4228.9
$ws.Range("N73").Value = -7948.5386
$ws.Range("H114").Value = 45000
$ws.Range("J114").Value = 45000
$ws.Range("L114").Value = 45000
$ws.Range("N114").Value = -53678
$ws.Range("H122").Value = 7981.8887
$ws.Range("I122").Value = 17337.666
$ws.Range("K122").Value = 52012.99800000001
$ws.Range("M122").Value = -49562.99800000001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 853.375
$ws.Range("I22").Value = 901.6667
$ws.Range("K22").Value = 901.6667
$ws.Range("M22").Value = -606.6667
$ws.Range("H27").Value = 853.375
$ws.Range("I27").Value = 901.6667
$ws.Range("K27").Value = 901.6667
$ws.Range("M27").Value = -794.6667
$ws.Range("H40").Value = 3800
$ws.Range("I40").Value = 3450
$ws.Range("K40").Value = 3450
$ws.Range("M40").Value = -3314
$ws.Range("H133").Value = 64318.5
$ws.Range("J133").Value = 64318.5
$ws.Range("L133").Value = 64318.5
$ws.Range("N133").Value = -69378.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 18624.428
$ws.Range("I49").Value = 9999
$ws.Range("J49").Value = 20062
$ws.Range("K49").Value = 9999
$ws.Range("L49").Value = 20062
$ws.Range("M49").Value = -9769
$ws.Range("N49").Value = -20522
$ws.Range("H70").Value = 29773.4
$ws.Range("J70").Value = 29773.4
$ws.Range("L70").Value = 29773.4
$ws.Range("N70").Value = -30403.4
$ws.Range("H73").Value = 29773.4
$ws.Range("J73").Value = 29773.4
$ws.Range("L73").Value = 29773.4
$ws.Range("N73").Value = -31957.4

